# Fixed sync of different markets
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs")

# Update tickers (column B)
$ws.Range("B6").Value = "SAN"
$ws.Range("B5").Value = "REPYY"
$ws.Range("B4").Value = "SLR.MC"
$ws.Range("B7").Value = "MT"

# Update weights (column C)
$ws.Range("C4").Value = 0.3
$ws.Range("C5").Value = 0.3
$ws.Range("C6").Value = 0.25
$ws.Range("C7").Value = 0.15

# Update lower bounds (column D)
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 0.2
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 0.1

# Update upper bounds (column E)
$ws.Range("E4").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("E7").Value = 1

# Update start date (H4) -> 2017-03-01
$ws.Range("H4").Value = "2017-03-01"

# Page setup (paper size A4->Letter-style "9" = A4, portrait orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Re-select current active cell (H17) for sheet view
$ws.Range("H17").Select() | Out-Null
